$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header labels:
#   in-dev       -> open
#   completed    -> closed
#   % in-dev     -> % open
#   % completed  -> % closed
$ws.Range("D1").Value = "open"
$ws.Range("E1").Value = "closed"
$ws.Range("G1").Value = "% open"
$ws.Range("H1").Value = "% closed"
